$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2564746666666666
$ws.Range("H2").Value = 0.7694239999999999
$ws.Range("I2").Value = 0.1818007399394835
$ws.Range("J2").Value = 0.1818007399394835
$ws.Range("M2").Value = 3.759736666666667
$ws.Range("N2").Value = 11.27921
$ws.Range("O2").Value = 0.0683751702595819
$ws.Range("P2").Value = 0.06837517025958188
$ws.Range("Q2").Value = 0.9642772083377776
$ws.Range("R2").Value = 8.678494875039998
$ws.Range("S2").Value = 0.01243065654668016
$ws.Range("T2").Value = 0.01243065654668015

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2564746666666666
$ws.Range("H3").Value = 0.7694239999999999
$ws.Range("I3").Value = 0.1818007399394835
$ws.Range("J3").Value = 0.1818007399394835
$ws.Range("O3").Value = 0.6514180024294648
$ws.Range("P3").Value = 0.6514180024294647
$ws.Range("Q3").Value = 9.186778335745776
$ws.Range("R3").Value = 82.681005021712
$ws.Range("S3").Value = 0.118428274851577
$ws.Range("T3").Value = 0.118428274851577

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2564746666666666
$ws.Range("H4").Value = 0.7694239999999999
$ws.Range("I4").Value = 0.1818007399394835
$ws.Range("J4").Value = 0.1818007399394835
$ws.Range("O4").Value = 0.2802068273109533
$ws.Range("P4").Value = 0.2802068273109533
$ws.Range("Q4").Value = 3.951683866684444
$ws.Range("R4").Value = 35.56515480015999
$ws.Range("S4").Value = 0.05094180854122639
$ws.Range("T4").Value = 0.05094180854122639

$ws.Range("I5").Value = 0.7694380609030022
$ws.Range("J5").Value = 0.7694380609030022
$ws.Range("M5").Value = 3.759736666666667
$ws.Range("N5").Value = 11.27921
$ws.Range("O5").Value = 0.0683751702595819
$ws.Range("P5").Value = 0.06837517025958188
$ws.Range("Q5").Value = 4.081125223161111
$ws.Range("R5").Value = 36.73012700845
$ws.Range("S5").Value = 0.05261045841844532
$ws.Range("T5").Value = 0.05261045841844531

$ws.Range("I6").Value = 0.7694380609030022
$ws.Range("J6").Value = 0.7694380609030022
$ws.Range("O6").Value = 0.6514180024294648
$ws.Range("P6").Value = 0.6514180024294647
$ws.Range("S6").Value = 0.5012258046266346
$ws.Range("T6").Value = 0.5012258046266345

$ws.Range("I7").Value = 0.7694380609030022
$ws.Range("J7").Value = 0.7694380609030022
$ws.Range("O7").Value = 0.2802068273109533
$ws.Range("P7").Value = 0.2802068273109533
$ws.Range("S7").Value = 0.2156017978579223
$ws.Range("T7").Value = 0.2156017978579222

$ws.Range("G8").Value = 0.06878966666666667
$ws.Range("I8").Value = 0.0487611991575143
$ws.Range("J8").Value = 0.0487611991575143
$ws.Range("M8").Value = 3.759736666666667
$ws.Range("N8").Value = 11.27921
$ws.Range("O8").Value = 0.0683751702595819
$ws.Range("P8").Value = 0.06837517025958188
$ws.Range("Q8").Value = 0.2586310320544444
$ws.Range("R8").Value = 2.32767928849
$ws.Range("S8").Value = 0.003334055294456421
$ws.Range("T8").Value = 0.003334055294456421

$ws.Range("G9").Value = 0.06878966666666667
$ws.Range("I9").Value = 0.0487611991575143
$ws.Range("J9").Value = 0.0487611991575143
$ws.Range("O9").Value = 0.6514180024294648
$ws.Range("P9").Value = 0.6514180024294647
$ws.Range("S9").Value = 0.03176392295125326
$ws.Range("T9").Value = 0.03176392295125326

$ws.Range("G10").Value = 0.06878966666666667
$ws.Range("I10").Value = 0.0487611991575143
$ws.Range("J10").Value = 0.0487611991575143
$ws.Range("O10").Value = 0.2802068273109533
$ws.Range("P10").Value = 0.2802068273109533
$ws.Range("S10").Value = 0.01366322091180461
$ws.Range("T10").Value = 0.01366322091180461

Write-Host "Applied TPM update"